$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'277.47"
$ws.Range("E2").Value = "'0.81%"
$ws.Range("D3").Value = "'27.31"
$ws.Range("E3").Value = "'2.40%"
$ws.Range("D4").Value = "'4.892"
$ws.Range("E4").Value = "'0.73%"
$ws.Range("E5").Value = "'1.11%"
$ws.Range("D6").Value = "'6.963"
$ws.Range("E6").Value = "'1.13%"
$ws.Range("D7").Value = "'0.8835"
$ws.Range("E7").Value = "'1.84%"
$ws.Range("E8").Value = "'-6.92%"
$ws.Range("D9").Value = "'0.1538"
$ws.Range("E9").Value = "'-1.63%"
$ws.Range("D10").Value = "'0.05140"
$ws.Range("E10").Value = "'-1.38%"
$ws.Range("D11").Value = "'0.07422"
$ws.Range("E11").Value = "'0.01%"
$ws.Range("D12").Value = "'0.02898"
$ws.Range("E12").Value = "'-1.11%"
$ws.Range("D13").Value = "'0.08972"
$ws.Range("E13").Value = "'-0.60%"
$ws.Range("D14").Value = "'0.001565"
$ws.Range("E14").Value = "'-0.41%"
$ws.Range("D15").Value = "'0.0006368"
$ws.Range("E15").Value = "'0.88%"
$ws.Range("D16").Value = "'0.006128"
$ws.Range("E16").Value = "'2.74%"
$ws.Range("D17").Value = "'3.477"
$ws.Range("E17").Value = "'0.79%"
$ws.Range("D18").Value = "'3.309"
$ws.Range("E18").Value = "'-0.10%"
$ws.Range("D19").Value = "'2.273"
$ws.Range("E19").Value = "'0.04%"
$ws.Range("D20").Value = "'0.3152"
$ws.Range("E20").Value = "'1.23%"
$ws.Range("E21").Value = "'1.11%"
$ws.Range("D22").Value = "'3.912"
$ws.Range("E22").Value = "'0.34%"
$ws.Range("D23").Value = "'0.04413"
$ws.Range("E23").Value = "'1.50%"
$ws.Range("D24").Value = "'0.1501"
$ws.Range("E24").Value = "'8.75%"
$ws.Range("E26").Value = "'-0.01%"
$ws.Range("E27").Value = "'-9.01%"
$ws.Range("D28").Value = "'0.0001181"
$ws.Range("E28").Value = "'-1.55%"
$ws.Range("E29").Value = "'15.68%"
$ws.Range("D40").Value = "'0.04148"
$ws.Range("E40").Value = "'0.72%"
$ws.Range("D41").Value = "'0.006795"
$ws.Range("E41").Value = "'-0.14%"
$ws.Range("D42").Value = "'0.1174"
$ws.Range("E42").Value = "'0.48%"
$ws.Range("D43").Value = "'0.002011"
$ws.Range("E43").Value = "'-6.18%"
$ws.Range("D44").Value = "'0.01149"
$ws.Range("E44").Value = "'7.02%"
$ws.Range("D45").Value = "'0.00005318"
$ws.Range("E45").Value = "'0.64%"
$ws.Range("D46").Value = "'1.687"
$ws.Range("E46").Value = "'13.27%"
$ws.Range("D47").Value = "'0.01853"
$ws.Range("E47").Value = "'-11.71%"
